$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "T7pol_" prefix from the variant names in column A (rows 2-10)
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value2
    $cell.Value = $cur.Replace("T7pol_", "")
}

# Update header for column B
$ws.Range("B1").Value = "fitness"

# Move the active selection to B2
$ws.Range("B2").Select()
